$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two obsolete items (rows 8 and 9: "SYRUP 120 ML" and "EXTRA TAB").
# Everything below shifts up by two rows.
$ws.Rows("8:9").Delete()

# Remove the obsolete "سرنجات 5 سم" item, which (after the shift above) now
# lives on row 12.
$ws.Rows("12").Delete()

# Row 7 ("1 2 3 (ONE TWO THREE) 20 F.C.TABS.") balance/price/count updates.
$ws.Range("H7").Value = "10:0"
$ws.Range("Q7").Value = "0:-1"

# P7 is formatted with a numeric (0.00) display format, so a plain string
# assignment of "-20.0000" would be auto-coerced to the number -20. Force it
# to stay text (matching the source data, which stores it as text) without
# disturbing the cell's existing number format/style.
$p7Format = $ws.Range("P7").NumberFormat
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "-20.0000"
$ws.Range("P7").NumberFormat = $p7Format

# Renumber the remaining rows sequentially (column A item index).
$ws.Range("A8").Value = 2
$ws.Range("A9").Value = 3
$ws.Range("A10").Value = 4
$ws.Range("A11").Value = 5

# Row heights follow an alternating pattern (odd=25.5 / even=24.75); restore
# it for the two rows whose heights came along with the deleted rows above.
$ws.Rows("10").RowHeight = 24.75
$ws.Rows("11").RowHeight = 25.5

# Update the grand-total row (now row 12) for the removed items' prices and
# give it its own distinct row height.
$ws.Range("N12").Value = 157.14
$ws.Rows("12").RowHeight = 26.25
